# Apply updated TPM values to Fgf2-Fgfr2 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2858606666666667
$ws.Range("N2").Value = 0.857582
$ws.Range("O2").Value = 0.0687156860066334
$ws.Range("P2").Value = 0.06932858672617494
$ws.Range("Q2").Value = 0.1351371998386667
$ws.Range("R2").Value = 1.216234798548
$ws.Range("S2").Value = 0.002250884653560849
$ws.Range("T2").Value = 0.00235518832244885
$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.9046431256549901
$ws.Range("P3").Value = 0.9127119736118995
$ws.Range("Q3").Value = 1.779083437259334
$ws.Range("R3").Value = 16.011750935334
$ws.Range("S3").Value = 0.02963293313101129
$ws.Range("T3").Value = 0.03100609263102728
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110331
$ws.Range("N4").Value = 0.220662
$ws.Range("O4").Value = 0.02652155835639462
$ws.Range("P4").Value = 0.01783874265571248
$ws.Range("Q4").Value = 0.052157656278
$ws.Range("R4").Value = 0.312945937668
$ws.Range("S4").Value = 0.0008687531502947427
$ws.Range("T4").Value = 0.0006060068490339212
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.0004976666666666667
$ws.Range("N5").Value = 0.001493
$ws.Range("O5").Value = 0.0001196299819817856
$ws.Range("P5").Value = 0.0001206970062130259
$ws.Range("Q5").Value = 0.0002352659446666667
$ws.Range("R5").Value = 0.002117393502
$ws.Range("S5").Value = 0.000003918658259812295
$ws.Range("T5").Value = 0.000004100244834215426
$ws.Range("I6").Value = 0.822180234441485
$ws.Range("J6").Value = 0.8526739017519405
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2858606666666667
$ws.Range("N6").Value = 0.857582
$ws.Range("O6").Value = 0.0687156860066334
$ws.Range("P6").Value = 0.06932858672617494
$ws.Range("Q6").Value = 3.391912137875555
$ws.Range("R6").Value = 30.52720924088
$ws.Range("S6").Value = 0.05649667883074132
$ws.Range("T6").Value = 0.05911467654675538
$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("O7").Value = 0.9046431256549901
$ws.Range("P7").Value = 0.9127119736118995
$ws.Range("S7").Value = 0.7437796971368975
$ws.Range("T7").Value = 0.7782456797153725
$ws.Range("D8").Value = "MuSCs"
$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.110331
$ws.Range("N8").Value = 0.220662
$ws.Range("O8").Value = 0.02652155835639462
$ws.Range("P8").Value = 0.01783874265571248
$ws.Range("Q8").Value = 1.30914498468
$ws.Range("R8").Value = 7.85486990808
$ws.Range("S8").Value = 0.02180550106721405
$ws.Range("T8").Value = 0.01521063030259513
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.0004976666666666667
$ws.Range("N9").Value = 0.001493
$ws.Range("O9").Value = 0.0001196299819817856
$ws.Range("P9").Value = 0.0001206970062130259
$ws.Range("Q9").Value = 0.005905120235555556
$ws.Range("R9").Value = 0.05314608212
$ws.Range("S9").Value = 0.00009835740663201512
$ws.Range("T9").Value = 0.000102915187217439
$ws.Range("G10").Value = 0.37892
$ws.Range("H10").Value = 1.13676
$ws.Range("I10").Value = 0.02625574638939025
$ws.Range("J10").Value = 0.02722954016579943
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2858606666666667
$ws.Range("N10").Value = 0.857582
$ws.Range("O10").Value = 0.0687156860066334
$ws.Range("P10").Value = 0.06932858672617494
$ws.Range("Q10").Value = 0.1083183238133333
$ws.Range("R10").Value = 0.97486491432
$ws.Range("S10").Value = 0.001804181624763139
$ws.Range("T10").Value = 0.00188778553689849
$ws.Range("G11").Value = 0.37892
$ws.Range("H11").Value = 1.13676
$ws.Range("I11").Value = 0.02625574638939025
$ws.Range("J11").Value = 0.02722954016579943
$ws.Range("O11").Value = 0.9046431256549901
$ws.Range("P11").Value = 0.9127119736118995
$ws.Range("Q11").Value = 1.426012497506667
$ws.Range("R11").Value = 12.83411247756
$ws.Range("S11").Value = 0.02375208048010272
$ws.Range("T11").Value = 0.02485272734527129
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.110331
$ws.Range("N12").Value = 0.220662
$ws.Range("O12").Value = 0.02652155835639462
$ws.Range("P12").Value = 0.01783874265571248
$ws.Range("Q12").Value = 0.04180662252
$ws.Range("R12").Value = 0.25083973512
$ws.Range("S12").Value = 0.0006963433100569107
$ws.Range("T12").Value = 0.0004857407596510824
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.0004976666666666667
$ws.Range("N13").Value = 0.001493
$ws.Range("O13").Value = 0.0001196299819817856
$ws.Range("P13").Value = 0.0001206970062130259
$ws.Range("Q13").Value = 0.0001885758533333333
$ws.Range("R13").Value = 0.00169718268
$ws.Range("S13").Value = 0.000003140974467481088
$ws.Range("T13").Value = 0.000003286523978569332
$ws.Range("G14").Value = 1.548357
$ws.Range("H14").Value = 3.096714
$ws.Range("I14").Value = 0.1072872076222874
$ws.Range("J14").Value = 0.0741775733180209
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.2858606666666667
$ws.Range("N14").Value = 0.857582
$ws.Range("O14").Value = 0.0687156860066334
$ws.Range("P14").Value = 0.06932858672617494
$ws.Range("Q14").Value = 0.442614364258
$ws.Range("R14").Value = 2.655686185548
$ws.Range("S14").Value = 0.007372314071501584
$ws.Range("T14").Value = 0.005142626324915612
$ws.Range("G15").Value = 1.548357
$ws.Range("H15").Value = 3.096714
$ws.Range("I15").Value = 0.1072872076222874
$ws.Range("J15").Value = 0.0741775733180209
$ws.Range("O15").Value = 0.9046431256549901
$ws.Range("P15").Value = 0.9127119736118995
$ws.Range("Q15").Value = 5.827025315639001
$ws.Range("R15").Value = 34.962151893834
$ws.Range("S15").Value = 0.0970566348462219
$ws.Range("T15").Value = 0.06770275934083222
$ws.Range("D16").Value = "MuSCs"
$ws.Range("G16").Value = 1.548357
$ws.Range("H16").Value = 3.096714
$ws.Range("I16").Value = 0.1072872076222874
$ws.Range("J16").Value = 0.0741775733180209
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.110331
$ws.Range("N16").Value = 0.220662
$ws.Range("O16").Value = 0.02652155835639462
$ws.Range("P16").Value = 0.01783874265571248
$ws.Range("Q16").Value = 0.170831776167
$ws.Range("R16").Value = 0.6833271046680001
$ws.Range("S16").Value = 0.002845423937849119
$ws.Range("T16").Value = 0.001323234641245419
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.0004976666666666667
$ws.Range("N17").Value = 0.001493
$ws.Range("O17").Value = 0.0001196299819817856
$ws.Range("P17").Value = 0.0001206970062130259
$ws.Range("Q17").Value = 0.0007705656670000001
$ws.Range("R17").Value = 0.004623394002
$ws.Range("S17").Value = 0.00001283476671473033
$ws.Range("T17").Value = 0.000008953011027632353
$ws.Range("G18").Value = 0.16626
$ws.Range("H18").Value = 0.49878
$ws.Range("I18").Value = 0.01152032195371061
$ws.Range("J18").Value = 0.01194759671689489
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.2858606666666667
$ws.Range("N18").Value = 0.857582
$ws.Range("O18").Value = 0.0687156860066334
$ws.Range("P18").Value = 0.06932858672617494
$ws.Range("Q18").Value = 0.04752719444
$ws.Range("R18").Value = 0.42774474996
$ws.Range("S18").Value = 0.0007916268260665035
$ws.Range("T18").Value = 0.0008283099951566106
$ws.Range("G19").Value = 0.16626
$ws.Range("H19").Value = 0.49878
$ws.Range("I19").Value = 0.01152032195371061
$ws.Range("J19").Value = 0.01194759671689489
$ws.Range("O19").Value = 0.9046431256549901
$ws.Range("P19").Value = 0.9127119736118995
$ws.Range("Q19").Value = 0.62569628902
$ws.Range("R19").Value = 5.63126660118
$ws.Range("S19").Value = 0.01042178006075656
$ws.Range("T19").Value = 0.01090471457939619
$ws.Range("D20").Value = "MuSCs"
$ws.Range("G20").Value = 0.16626
$ws.Range("H20").Value = 0.49878
$ws.Range("I20").Value = 0.01152032195371061
$ws.Range("J20").Value = 0.01194759671689489
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.110331
$ws.Range("N20").Value = 0.220662
$ws.Range("O20").Value = 0.02652155835639462
$ws.Range("P20").Value = 0.01783874265571248
$ws.Range("Q20").Value = 0.01834363206
$ws.Range("R20").Value = 0.11006179236
$ws.Range("S20").Value = 0.0003055368909797899
$ws.Range("T20").Value = 0.0002131301031869233
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("G21").Value = 0.16626
$ws.Range("H21").Value = 0.49878
$ws.Range("I21").Value = 0.01152032195371061
$ws.Range("J21").Value = 0.01194759671689489
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.0004976666666666667
$ws.Range("N21").Value = 0.001493
$ws.Range("O21").Value = 0.0001196299819817856
$ws.Range("P21").Value = 0.0001206970062130259
$ws.Range("Q21").Value = 0.00008274205999999999
$ws.Range("R21").Value = 0.00074467854
$ws.Range("S21").Value = 0.000001378175907746769
$ws.Range("T21").Value = 0.000001442039155169791
